$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 5000
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 5000
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 5000
$ws.Cells.Item(40, 14).Value = -5350
$ws.Cells.Item(40, 13).ClearContents()

$ws.Cells.Item(111, 8).Value = 3602.25
$ws.Cells.Item(111, 10).Value = 2599.8
$ws.Cells.Item(111, 12).Value = 7799.400000000001
$ws.Cells.Item(111, 14).Value = -13933.4

$ws.Cells.Item(112, 8).Value = 3623.077
$ws.Cells.Item(112, 10).Value = 3833.3333
$ws.Cells.Item(112, 12).Value = 11499.9999
$ws.Cells.Item(112, 14).Value = -13715.9999

$ws.Cells.Item(132, 8).Value = 7938580
$ws.Cells.Item(132, 9).Value = 8132052.5
$ws.Cells.Item(132, 11).Value = 24396157.5
$ws.Cells.Item(132, 13).Value = -24393627.5

$ws.Cells.Item(137, 8).Value = 1380.8695
$ws.Cells.Item(137, 9).Value = 1290.75
$ws.Cells.Item(137, 10).Value = 1586.8572
$ws.Cells.Item(137, 11).Value = 3872.25
$ws.Cells.Item(137, 12).Value = 4760.571599999999
$ws.Cells.Item(137, 13).Value = -1322.25
$ws.Cells.Item(137, 14).Value = -9860.571599999999

$ws.Cells.Item(140, 8).Value = 30710
$ws.Cells.Item(140, 10).Value = 30710
$ws.Cells.Item(140, 12).Value = 30710
$ws.Cells.Item(140, 14).Value = -41070

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 26052.75

$ws.Cells.Item(116, 8).Value = 26052.75

$ws.Cells.Item(119, 8).Value = 24499.5
$ws.Cells.Item(119, 10).Value = 24499.5
$ws.Cells.Item(119, 12).Value = 24499.5
$ws.Cells.Item(119, 14).Value = -34175.5

$ws.Cells.Item(132, 8).Value = 2628.6538
$ws.Cells.Item(132, 9).Value = 2466.7144
$ws.Cells.Item(132, 11).Value = 7400.1432
$ws.Cells.Item(132, 13).Value = -4870.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 26052.75

$ws.Cells.Item(82, 8).Value = 29451
$ws.Cells.Item(82, 10).Value = 33500
$ws.Cells.Item(82, 12).Value = 33500
$ws.Cells.Item(82, 14).Value = -34266

$ws.Cells.Item(85, 8).Value = 29451
$ws.Cells.Item(85, 10).Value = 33500
$ws.Cells.Item(85, 12).Value = 33500
$ws.Cells.Item(85, 14).Value = -36152

$ws.Cells.Item(99, 8).Value = 62501588
$ws.Cells.Item(99, 9).Value = 83334830
$ws.Cells.Item(99, 10).Value = 1850
$ws.Cells.Item(99, 11).Value = 83334830
$ws.Cells.Item(99, 12).Value = 1850
$ws.Cells.Item(99, 13).Value = -83333332
$ws.Cells.Item(99, 14).Value = -4846

$ws.Cells.Item(107, 8).Value = 1652.3077
$ws.Cells.Item(107, 9).Value = 995.7143
$ws.Cells.Item(107, 11).Value = 995.7143
$ws.Cells.Item(107, 13).Value = 924.2857

$ws.Cells.Item(134, 8).Value = 8480.294
$ws.Cells.Item(134, 9).Value = 930.5
$ws.Cells.Item(134, 10).Value = 26599.8
$ws.Cells.Item(134, 11).Value = 2791.5
$ws.Cells.Item(134, 12).Value = 79799.39999999999
$ws.Cells.Item(134, 13).Value = -256.5
$ws.Cells.Item(134, 14).Value = -84869.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1274.75
$ws.Cells.Item(31, 9).Value = 1274.75
$ws.Cells.Item(31, 10).Value = 0
$ws.Cells.Item(31, 11).Value = 1274.75
$ws.Cells.Item(31, 12).Value = 0
$ws.Cells.Item(31, 13).Value = -979.75
$ws.Cells.Item(31, 14).ClearContents()

$ws.Cells.Item(34, 8).Value = 1274.75
$ws.Cells.Item(34, 9).Value = 1274.75
$ws.Cells.Item(34, 10).Value = 0
$ws.Cells.Item(34, 11).Value = 1274.75
$ws.Cells.Item(34, 12).Value = 0
$ws.Cells.Item(34, 13).Value = -1072.75
$ws.Cells.Item(34, 14).ClearContents()

$ws.Cells.Item(58, 8).Value = 949.5714
$ws.Cells.Item(58, 9).Value = 703.0357
$ws.Cells.Item(58, 10).Value = 1935.7142
$ws.Cells.Item(58, 11).Value = 703.0357
$ws.Cells.Item(58, 12).Value = 1935.7142
$ws.Cells.Item(58, 13).Value = -500.0357
$ws.Cells.Item(58, 14).Value = -2341.7142

$ws.Cells.Item(107, 8).Value = 607.6
$ws.Cells.Item(107, 9).Value = 467.53333
$ws.Cells.Item(107, 10).Value = 817.7
$ws.Cells.Item(107, 11).Value = 467.53333
$ws.Cells.Item(107, 12).Value = 817.7
$ws.Cells.Item(107, 13).Value = 1452.46667
$ws.Cells.Item(107, 14).Value = -4657.7

$ws.Cells.Item(132, 8).Value = 6026.1377
$ws.Cells.Item(132, 9).Value = 9460.23
$ws.Cells.Item(132, 10).Value = 3235.9375
$ws.Cells.Item(132, 11).Value = 28380.69
$ws.Cells.Item(132, 12).Value = 9707.8125
$ws.Cells.Item(132, 13).Value = -25850.69
$ws.Cells.Item(132, 14).Value = -14767.8125

$ws.Cells.Item(136, 8).Value = 949.5714
$ws.Cells.Item(136, 9).Value = 703.0357
$ws.Cells.Item(136, 10).Value = 1935.7142
$ws.Cells.Item(136, 11).Value = 2109.1071
$ws.Cells.Item(136, 12).Value = 5807.142599999999
$ws.Cells.Item(136, 13).Value = 440.8928999999998
$ws.Cells.Item(136, 14).Value = -10907.1426

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 164.5
$ws.Cells.Item(2, 10).Value = 211.75
$ws.Cells.Item(2, 12).Value = 1270.5
$ws.Cells.Item(2, 14).Value = -1496.5

$ws.Cells.Item(6, 8).Value = 836.75
$ws.Cells.Item(6, 9).Value = 173.5
$ws.Cells.Item(6, 11).Value = 520.5
$ws.Cells.Item(6, 13).Value = -407.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 15000

$ws.Cells.Item(81, 8).Value = 15000

$ws.Cells.Item(84, 8).Value = 15000

$ws.Cells.Item(93, 8).Value = 29999.8
$ws.Cells.Item(93, 10).Value = 29999.8
$ws.Cells.Item(93, 12).Value = 29999.8
$ws.Cells.Item(93, 14).Value = -33743.8

$ws.Cells.Item(96, 8).Value = 30261
$ws.Cells.Item(96, 10).Value = 30261
$ws.Cells.Item(96, 12).Value = 30261
$ws.Cells.Item(96, 14).Value = -35753

$ws.Cells.Item(122, 8).Value = 302520
$ws.Cells.Item(122, 9).Value = 3700
$ws.Cells.Item(122, 10).Value = 750750
$ws.Cells.Item(122, 11).Value = 11100
$ws.Cells.Item(122, 12).Value = 2252250
$ws.Cells.Item(122, 13).Value = -8650
$ws.Cells.Item(122, 14).Value = -2257150

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3500
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 13).ClearContents()

$ws.Cells.Item(126, 8).Value = 3500
$ws.Cells.Item(126, 9).Value = 0
$ws.Cells.Item(126, 11).Value = 0
$ws.Cells.Item(126, 13).ClearContents()

$ws.Cells.Item(136, 8).Value = 5685.913
$ws.Cells.Item(136, 9).Value = 7267.4375
$ws.Cells.Item(136, 10).Value = 2071
$ws.Cells.Item(136, 11).Value = 21802.3125
$ws.Cells.Item(136, 12).Value = 6213
$ws.Cells.Item(136, 13).Value = -19252.3125
$ws.Cells.Item(136, 14).Value = -11313

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 660
$ws.Cells.Item(81, 9).Value = 660
$ws.Cells.Item(81, 11).Value = 1320
$ws.Cells.Item(81, 13).Value = -259

$ws.Cells.Item(84, 8).Value = 660
$ws.Cells.Item(84, 9).Value = 660
$ws.Cells.Item(84, 11).Value = 6600
$ws.Cells.Item(84, 13).Value = -1296

$ws.Cells.Item(92, 8).Value = 27500
$ws.Cells.Item(92, 10).Value = 27500
$ws.Cells.Item(92, 12).Value = 27500
$ws.Cells.Item(92, 14).Value = -32492

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()

$ws.Cells.Item(98, 8).Value = 0
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(99, 8).Value = 15500
$ws.Cells.Item(99, 10).Value = 15500
$ws.Cells.Item(99, 12).Value = 15500
$ws.Cells.Item(99, 14).Value = -21490

$ws.Cells.Item(101, 8).Value = 13366.833
$ws.Cells.Item(101, 10).Value = 13366.833
$ws.Cells.Item(101, 12).Value = 13366.833
$ws.Cells.Item(101, 14).Value = -19856.833

$ws.Cells.Item(105, 8).Value = 34998.5
$ws.Cells.Item(105, 10).Value = 34998.5
$ws.Cells.Item(105, 12).Value = 34998.5
$ws.Cells.Item(105, 14).Value = -41986.5

$ws.Cells.Item(123, 8).Value = 53200
$ws.Cells.Item(123, 10).Value = 53200
$ws.Cells.Item(123, 12).Value = 53200
$ws.Cells.Item(123, 14).Value = -63000
